$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header "Save" in H1, copying the header formatting (bold, border, alignment)
# from the existing "sum" header cell (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add 0 values in H2:H12 for the new "Save" column
$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
